# Generate Report for handoff
#
# The "1d3b8b86-2b3c-4139-a745-190078635726..." row (row 4) in both the
# zh-cn and de-de localization-status sheets just went through a fresh
# handoff generation, so its "Latest Handoff Datetime" cell (column D)
# is stamped with a new timestamp while the rest of the row (file names,
# status, etc.) stays as-is.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-20 06:45:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-20 06:45:26"
